$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new log row (row 8) by copying the formatting of the prior row (row 7)
# and then overwriting the values/text that differ:
#   Date   : same day (1/29/2018) as row 7
#   Resource : Drew (same shared string as other rows)
#   Time Spent (mins) : 60
#   Description : new entry about the ECS work done during Data Comm class
$ws.Range("A7:D7").Copy($ws.Range("A8:D8"))

$ws.Range("C8").Value = 60
$ws.Range("D8").Value = "Programming basic ECS during Data Communications Class"

# Leave the active selection on the newly added description cell
$ws.Range("D8").Select()

$wb.Save()
